$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 24.863113
$ws.Range("H2").Value = 49.726226
$ws.Range("I2").Value = 0.0840270138440103
$ws.Range("J2").Value = 0.05854365446293047
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 2.019046
$ws.Range("N2").Value = 4.038092
$ws.Range("O2").Value = 0.003975353327590414
$ws.Range("P2").Value = 0.002707057536954368
$ws.Range("Q2").Value = 50.19976885019799
$ws.Range("R2").Value = 200.799075400792
$ws.Range("S2").Value = 0.0003340370690922721
$ws.Range("T2").Value = 0.0001584810410547281

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 24.863113
$ws.Range("H3").Value = 49.726226
$ws.Range("I3").Value = 0.0840270138440103
$ws.Range("J3").Value = 0.05854365446293047
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 139.6948166666666
$ws.Range("N3").Value = 419.0844499999999
$ws.Range("O3").Value = 0.2750488370661026
$ws.Range("P3").Value = 0.2809459811695414
$ws.Range("Q3").Value = 3473.248012297616
$ws.Range("R3").Value = 20839.48807378569
$ws.Range("S3").Value = 0.02311153243993233
$ws.Range("T3").Value = 0.01644760444433861

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 24.863113
$ws.Range("H4").Value = 49.726226
$ws.Range("I4").Value = 0.0840270138440103
$ws.Range("J4").Value = 0.05854365446293047
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 186.3548536666667
$ws.Range("N4").Value = 559.064561
$ws.Range("O4").Value = 0.3669190239530987
$ws.Range("P4").Value = 0.3747858972750337
$ws.Range("Q4").Value = 4633.361784812798
$ws.Range("R4").Value = 27800.17070887679
$ws.Range("S4").Value = 0.03083110990533778
$ws.Range("T4").Value = 0.02194133606764893

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 24.863113
$ws.Range("H5").Value = 49.726226
$ws.Range("I5").Value = 0.0840270138440103
$ws.Range("J5").Value = 0.05854365446293047
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 143.6051993333333
$ws.Range("N5").Value = 430.815598
$ws.Range("O5").Value = 0.2827480934208787
$ws.Range("P5").Value = 0.2888103122968479
$ws.Range("Q5").Value = 3570.472298412191
$ws.Range("R5").Value = 21422.83379047315
$ws.Range("S5").Value = 0.02375847796024369
$ws.Range("T5").Value = 0.01690801112843771

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 24.863113
$ws.Range("H6").Value = 49.726226
$ws.Range("I6").Value = 0.0840270138440103
$ws.Range("J6").Value = 0.05854365446293047
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.253715333333335
$ws.Range("N6").Value = 18.761146
$ws.Range("O6").Value = 0.01231310631861279
$ws.Range("P6").Value = 0.0125771036621259
$ws.Range("Q6").Value = 155.4868310024994
$ws.Range("R6").Value = 932.9209860149962
$ws.Range("S6").Value = 0.001034633555096847
$ws.Range("T6").Value = 0.0007363096109399563

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 24.863113
$ws.Range("H7").Value = 49.726226
$ws.Range("I7").Value = 0.0840270138440103
$ws.Range("J7").Value = 0.05854365446293047
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 29.9633245
$ws.Range("N7").Value = 59.926649
$ws.Range("O7").Value = 0.05899558591371687
$ws.Range("P7").Value = 0.04017364805949665
$ws.Range("Q7").Value = 744.9815228991685
$ws.Range("R7").Value = 2979.926091596674
$ws.Range("S7").Value = 0.004957222914307387
$ws.Range("T7").Value = 0.002351912170510549

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 31.94572966666667
$ws.Range("H8").Value = 95.837189
$ws.Range("I8").Value = 0.1079633217673913
$ws.Range("J8").Value = 0.1128309893759997
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 2.019046
$ws.Range("N8").Value = 4.038092
$ws.Range("O8").Value = 0.003975353327590414
$ws.Range("P8").Value = 0.002707057536954368
$ws.Range("Q8").Value = 64.49989770056466
$ws.Range("R8").Value = 386.999386203388
$ws.Range("S8").Value = 0.0004291923504457137
$ws.Range("T8").Value = 0.0003054399801923183

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 31.94572966666667
$ws.Range("H9").Value = 95.837189
$ws.Range("I9").Value = 0.1079633217673913
$ws.Range("J9").Value = 0.1128309893759997
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 139.6948166666666
$ws.Range("N9").Value = 419.0844499999999
$ws.Range("O9").Value = 0.2750488370661026
$ws.Range("P9").Value = 0.2809459811695414
$ws.Range("Q9").Value = 4462.652849067893
$ws.Range("R9").Value = 40163.87564161104
$ws.Range("S9").Value = 0.02969518609791442
$ws.Range("T9").Value = 0.03169941301657034

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 31.94572966666667
$ws.Range("H10").Value = 95.837189
$ws.Range("I10").Value = 0.1079633217673913
$ws.Range("J10").Value = 0.1128309893759997
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 186.3548536666667
$ws.Range("N10").Value = 559.064561
$ws.Range("O10").Value = 0.3669190239530987
$ws.Range("P10").Value = 0.3747858972750337
$ws.Range("Q10").Value = 5953.241777306559
$ws.Range("R10").Value = 53579.17599575903
$ws.Range("S10").Value = 0.03961379664562557
$ws.Range("T10").Value = 0.04228746359371385

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 31.94572966666667
$ws.Range("H11").Value = 95.837189
$ws.Range("I11").Value = 0.1079633217673913
$ws.Range("J11").Value = 0.1128309893759997
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 143.6051993333333
$ws.Range("N11").Value = 430.815598
$ws.Range("O11").Value = 0.2827480934208787
$ws.Range("P11").Value = 0.2888103122968479
$ws.Range("Q11").Value = 4587.572876630446
$ws.Range("R11").Value = 41288.15588967402
$ws.Range("S11").Value = 0.03052642338911475
$ws.Range("T11").Value = 0.03258675327844481

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 31.94572966666667
$ws.Range("H12").Value = 95.837189
$ws.Range("I12").Value = 0.1079633217673913
$ws.Range("J12").Value = 0.1128309893759997
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.253715333333335
$ws.Range("N12").Value = 18.761146
$ws.Range("O12").Value = 0.01231310631861279
$ws.Range("P12").Value = 0.0125771036621259
$ws.Range("Q12").Value = 199.7794994509549
$ws.Range("R12").Value = 1798.015495058594
$ws.Range("S12").Value = 0.001329363859432492
$ws.Range("T12").Value = 0.001419087049682175

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 31.94572966666667
$ws.Range("H13").Value = 95.837189
$ws.Range("I13").Value = 0.1079633217673913
$ws.Range("J13").Value = 0.1128309893759997
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 29.9633245
$ws.Range("N13").Value = 59.926649
$ws.Range("O13").Value = 0.05899558591371687
$ws.Range("P13").Value = 0.04017364805949665
$ws.Range("Q13").Value = 957.2002643916101
$ws.Range("R13").Value = 5743.201586349661
$ws.Range("S13").Value = 0.006369359424858395
$ws.Range("T13").Value = 0.004532832457396218

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 95.81209199999999
$ws.Range("H14").Value = 287.436276
$ws.Range("I14").Value = 0.3238051478472381
$ws.Range("J14").Value = 0.3384043265671421
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 2.019046
$ws.Range("N14").Value = 4.038092
$ws.Range("O14").Value = 0.003975353327590414
$ws.Range("P14").Value = 0.002707057536954368
$ws.Range("Q14").Value = 193.449021104232
$ws.Range("R14").Value = 1160.694126625392
$ws.Range("S14").Value = 0.001287239871985424
$ws.Range("T14").Value = 0.0009160799827715493

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 95.81209199999999
$ws.Range("H15").Value = 287.436276
$ws.Range("I15").Value = 0.3238051478472381
$ws.Range("J15").Value = 0.3384043265671421
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 139.6948166666666
$ws.Range("N15").Value = 419.0844499999999
$ws.Range("O15").Value = 0.2750488370661026
$ws.Range("P15").Value = 0.2809459811695414
$ws.Range("Q15").Value = 13384.4526263898
$ws.Range("R15").Value = 120460.0736375082
$ws.Range("S15").Value = 0.08906222935140025
$ws.Range("T15").Value = 0.09507333555942366

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 95.81209199999999
$ws.Range("H16").Value = 287.436276
$ws.Range("I16").Value = 0.3238051478472381
$ws.Range("J16").Value = 0.3384043265671421
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 186.3548536666667
$ws.Range("N16").Value = 559.064561
$ws.Range("O16").Value = 0.3669190239530987
$ws.Range("P16").Value = 0.3747858972750337
$ws.Range("Q16").Value = 17855.04838415721
$ws.Range("R16").Value = 160695.4354574148
$ws.Range("S16").Value = 0.1188102687990974
$ws.Range("T16").Value = 0.1268291691742199

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 95.81209199999999
$ws.Range("H17").Value = 287.436276
$ws.Range("I17").Value = 0.3238051478472381
$ws.Range("J17").Value = 0.3384043265671421
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 143.6051993333333
$ws.Range("N17").Value = 430.815598
$ws.Range("O17").Value = 0.2827480934208787
$ws.Range("P17").Value = 0.2888103122968479
$ws.Range("Q17").Value = 13759.11457020367
$ws.Range("R17").Value = 123832.031131833
$ws.Range("S17").Value = 0.09155528819367233
$ws.Range("T17").Value = 0.09773465923846082

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 95.81209199999999
$ws.Range("H18").Value = 287.436276
$ws.Range("I18").Value = 0.3238051478472381
$ws.Range("J18").Value = 0.3384043265671421
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 6.253715333333335
$ws.Range("N18").Value = 18.761146
$ws.Range("O18").Value = 0.01231310631861279
$ws.Range("P18").Value = 0.0125771036621259
$ws.Range("Q18").Value = 599.1815488591441
$ws.Range("R18").Value = 5392.633939732296
$ws.Range("S18").Value = 0.003987047211957175
$ws.Range("T18").Value = 0.004256146294946853

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 95.81209199999999
$ws.Range("H19").Value = 287.436276
$ws.Range("I19").Value = 0.3238051478472381
$ws.Range("J19").Value = 0.3384043265671421
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 29.9633245
$ws.Range("N19").Value = 59.926649
$ws.Range("O19").Value = 0.05899558591371687
$ws.Range("P19").Value = 0.04017364805949665
$ws.Range("Q19").Value = 2870.848803619854
$ws.Range("R19").Value = 17225.09282171912
$ws.Range("S19").Value = 0.01910307441912553
$ws.Range("T19").Value = 0.01359493631731934

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 88.69318133333333
$ws.Range("H20").Value = 266.079544
$ws.Range("I20").Value = 0.2997461812511296
$ws.Range("J20").Value = 0.3132606299860782
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 2.019046
$ws.Range("N20").Value = 4.038092
$ws.Range("O20").Value = 0.003975353327590414
$ws.Range("P20").Value = 0.002707057536954368
$ws.Range("Q20").Value = 179.0756129983413
$ws.Range("R20").Value = 1074.453677990048
$ws.Range("S20").Value = 0.001191596979069198
$ws.Range("T20").Value = 0.0008480145494348866

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 88.69318133333333
$ws.Range("H21").Value = 266.079544
$ws.Range("I21").Value = 0.2997461812511296
$ws.Range("J21").Value = 0.3132606299860782
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 139.6948166666666
$ws.Range("N21").Value = 419.0844499999999
$ws.Range("O21").Value = 0.2750488370661026
$ws.Range("P21").Value = 0.2809459811695414
$ws.Range("Q21").Value = 12389.97770594342
$ws.Range("R21").Value = 111509.7993534908
$ws.Range("S21").Value = 0.0824448385681284
$ws.Range("T21").Value = 0.08800931505322743

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 88.69318133333333
$ws.Range("H22").Value = 266.079544
$ws.Range("I22").Value = 0.2997461812511296
$ws.Range("J22").Value = 0.3132606299860782
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 186.3548536666667
$ws.Range("N22").Value = 559.064561
$ws.Range("O22").Value = 0.3669190239530987
$ws.Range("P22").Value = 0.3747858972750337
$ws.Range("Q22").Value = 16528.40482860447
$ws.Range("R22").Value = 148755.6434574402
$ws.Range("S22").Value = 0.1099825762583331
$ws.Range("T22").Value = 0.1174056662902747

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 88.69318133333333
$ws.Range("H23").Value = 266.079544
$ws.Range("I23").Value = 0.2997461812511296
$ws.Range("J23").Value = 0.3132606299860782
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 143.6051993333333
$ws.Range("N23").Value = 430.815598
$ws.Range("O23").Value = 0.2827480934208787
$ws.Range("P23").Value = 0.2888103122968479
$ws.Range("Q23").Value = 12736.80198488081
$ws.Range("R23").Value = 114631.2178639273
$ws.Range("S23").Value = 0.08475266125894604
$ws.Range("T23").Value = 0.09047290037658659

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 88.69318133333333
$ws.Range("H24").Value = 266.079544
$ws.Range("I24").Value = 0.2997461812511296
$ws.Range("J24").Value = 0.3132606299860782
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 6.253715333333335
$ws.Range("N24").Value = 18.761146
$ws.Range("O24").Value = 0.01231310631861279
$ws.Range("P24").Value = 0.0125771036621259
$ws.Range("Q24").Value = 554.6619080663805
$ws.Range("R24").Value = 4991.957172597425
$ws.Range("S24").Value = 0.003690806598343338
$ws.Range("T24").Value = 0.003939911416597772

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 88.69318133333333
$ws.Range("H25").Value = 266.079544
$ws.Range("I25").Value = 0.2997461812511296
$ws.Range("J25").Value = 0.3132606299860782
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 29.9633245
$ws.Range("N25").Value = 59.926649
$ws.Range("O25").Value = 0.05899558591371687
$ws.Range("P25").Value = 0.04017364805949665
$ws.Range("Q25").Value = 2657.542573228009
$ws.Range("R25").Value = 15945.25543936805
$ws.Range("S25").Value = 0.01768370158830957
$ws.Range("T25").Value = 0.01258482229995691

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 41.14755033333334
$ws.Range("H26").Value = 123.442651
$ws.Range("I26").Value = 0.1390616606016355
$ws.Range("J26").Value = 0.1453314375020561
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 2.019046
$ws.Range("N26").Value = 4.038092
$ws.Range("O26").Value = 0.003975353327590414
$ws.Range("P26").Value = 0.002707057536954368
$ws.Range("Q26").Value = 83.07879691031533
$ws.Range("R26").Value = 498.472781461892
$ws.Range("S26").Value = 0.0005528192352129604
$ws.Range("T26").Value = 0.0003934205632463538

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 41.14755033333334
$ws.Range("H27").Value = 123.442651
$ws.Range("I27").Value = 0.1390616606016355
$ws.Range("J27").Value = 0.1453314375020561
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 139.6948166666666
$ws.Range("N27").Value = 419.0844499999999
$ws.Range("O27").Value = 0.2750488370661026
$ws.Range("P27").Value = 0.2809459811695414
$ws.Range("Q27").Value = 5748.099500097438
$ws.Range("R27").Value = 51732.89550087695
$ws.Range("S27").Value = 0.0382487480289609
$ws.Range("T27").Value = 0.04083028330379505

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 41.14755033333334
$ws.Range("H28").Value = 123.442651
$ws.Range("I28").Value = 0.1390616606016355
$ws.Range("J28").Value = 0.1453314375020561
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 186.3548536666667
$ws.Range("N28").Value = 559.064561
$ws.Range("O28").Value = 0.3669190239530987
$ws.Range("P28").Value = 0.3747858972750337
$ws.Range("Q28").Value = 7668.045721110136
$ws.Range("R28").Value = 69012.41148999122
$ws.Range("S28").Value = 0.05102436877724918
$ws.Range("T28").Value = 0.05446817320647859

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 41.14755033333334
$ws.Range("H29").Value = 123.442651
$ws.Range("I29").Value = 0.1390616606016355
$ws.Range("J29").Value = 0.1453314375020561
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 143.6051993333333
$ws.Range("N29").Value = 430.815598
$ws.Range("O29").Value = 0.2827480934208787
$ws.Range("P29").Value = 0.2888103122968479
$ws.Range("Q29").Value = 5909.0021676967
$ws.Range("R29").Value = 53181.01950927031
$ws.Range("S29").Value = 0.03931941940305376
$ws.Range("T29").Value = 0.04197321785151867

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 41.14755033333334
$ws.Range("H30").Value = 123.442651
$ws.Range("I30").Value = 0.1390616606016355
$ws.Range("J30").Value = 0.1453314375020561
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 6.253715333333335
$ws.Range("N30").Value = 18.761146
$ws.Range("O30").Value = 0.01231310631861279
$ws.Range("P30").Value = 0.0125771036621259
$ws.Range("Q30").Value = 257.3250664486719
$ws.Range("R30").Value = 2315.925598038047
$ws.Range("S30").Value = 0.001712281011830785
$ws.Range("T30").Value = 0.001827848554829132

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 41.14755033333334
$ws.Range("H31").Value = 123.442651
$ws.Range("I31").Value = 0.1390616606016355
$ws.Range("J31").Value = 0.1453314375020561
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 29.9633245
$ws.Range("N31").Value = 59.926649
$ws.Range("O31").Value = 0.05899558591371687
$ws.Range("P31").Value = 0.04017364805949665
$ws.Range("Q31").Value = 1232.91740301775
$ws.Range("R31").Value = 7397.5044181065
$ws.Range("S31").Value = 0.008204024145327922
$ws.Range("T31").Value = 0.005838494022188336

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 13.4326165
$ws.Range("H32").Value = 26.865233
$ws.Range("I32").Value = 0.04539667468859516
$ws.Range("J32").Value = 0.03162896210579336
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 2.019046
$ws.Range("N32").Value = 4.038092
$ws.Range("O32").Value = 0.003975353327590414
$ws.Range("P32").Value = 0.002707057536954368
$ws.Range("Q32").Value = 27.121070613859
$ws.Range("R32").Value = 108.484282455436
$ws.Range("S32").Value = 0.0001804678217848463
$ws.Range("T32").Value = 0.00008562142025453203

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 13.4326165
$ws.Range("H33").Value = 26.865233
$ws.Range("I33").Value = 0.04539667468859516
$ws.Range("J33").Value = 0.03162896210579336
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 139.6948166666666
$ws.Range("N33").Value = 419.0844499999999
$ws.Range("O33").Value = 0.2750488370661026
$ws.Range("P33").Value = 0.2809459811695414
$ws.Range("Q33").Value = 1876.466899321141
$ws.Range("R33").Value = 11258.80139592685
$ws.Range("S33").Value = 0.01248630257976627
$ws.Range("T33").Value = 0.008886029792186362

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 13.4326165
$ws.Range("H34").Value = 26.865233
$ws.Range("I34").Value = 0.04539667468859516
$ws.Range("J34").Value = 0.03162896210579336
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 186.3548536666667
$ws.Range("N34").Value = 559.064561
$ws.Range("O34").Value = 0.3669190239530987
$ws.Range("P34").Value = 0.3747858972750337
$ws.Range("Q34").Value = 2503.233282217952
$ws.Range("R34").Value = 15019.39969330771
$ws.Range("S34").Value = 0.01665690356745568
$ws.Range("T34").Value = 0.01185408894269781

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 13.4326165
$ws.Range("H35").Value = 26.865233
$ws.Range("I35").Value = 0.04539667468859516
$ws.Range("J35").Value = 0.03162896210579336
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 143.6051993333333
$ws.Range("N35").Value = 430.815598
$ws.Range("O35").Value = 0.2827480934208787
$ws.Range("P35").Value = 0.2888103122968479
$ws.Range("Q35").Value = 1928.993570050722
$ws.Range("R35").Value = 11573.96142030433
$ws.Range("S35").Value = 0.01283582321584814
$ws.Range("T35").Value = 0.00913477042339935

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 13.4326165
$ws.Range("H36").Value = 26.865233
$ws.Range("I36").Value = 0.04539667468859516
$ws.Range("J36").Value = 0.03162896210579336
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 6.253715333333335
$ws.Range("N36").Value = 18.761146
$ws.Range("O36").Value = 0.01231310631861279
$ws.Range("P36").Value = 0.0125771036621259
$ws.Range("Q36").Value = 84.00375977283635
$ws.Range("R36").Value = 504.0225586370181
$ws.Range("S36").Value = 0.0005589740819521503
$ws.Range("T36").Value = 0.0003978007351300151

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 13.4326165
$ws.Range("H37").Value = 26.865233
$ws.Range("I37").Value = 0.04539667468859516
$ws.Range("J37").Value = 0.03162896210579336
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 29.9633245
$ws.Range("N37").Value = 59.926649
$ws.Range("O37").Value = 0.05899558591371687
$ws.Range("P37").Value = 0.04017364805949665
$ws.Range("Q37").Value = 402.4858470735542
$ws.Range("R37").Value = 1609.943388294217
$ws.Range("S37").Value = 0.002678203421788072
$ws.Range("T37").Value = 0.001270650792125299

